# Update Moving File To Outside Directory
# - Change the Username id embedded in the free-text note cell (F2) from 33599 to 44912
# - Change the numeric Username cell (G2) from 33599 to 44912
# - Move the active selection from A2 to F2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the long descriptive text cell (F2) with the new username value
$ws.Range("F2").Value = "Username : 44912,`nPassword : bni1234,`nCetak Laporan PDF,`nNama Laporan : Aruskas ,`nTipe Laporan : Mutasi,`nProduk : - ,`nMata Uang : IDR,`nStatus Posting : Posting ,`nTanggal Transaksi : 01/08/2022,`nTanggal Pembanding : 31/07/2022"

# Update the numeric username cell (G2)
$ws.Range("G2").Value = 44912

# Update the selected cell in the sheet view from A2 to F2
$ws.Range("F2").Select()
